# Update "想去人数" (F column) values in both the "展览" and "全部类型"
# sheets, which contain duplicate data tables.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    2  = 70
    3  = 1057
    4  = 74
    5  = 3048
    7  = 2155
    10 = 1009
    13 = 248
    14 = 90
    15 = 95
    16 = 38
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
